$d = $word.ActiveDocument

# Locate the "Full-Stack Development and Data Engineering" paragraph under the
# Siege Analytics / PARTNER role so the new bullet points can be inserted
# immediately after it (and before the existing "Lead comprehensive research..." bullet).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -ceq "Full-Stack Development and Data Engineering`r") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Full-Stack Development and Data Engineering' paragraph."
}

$newBullets = @(
    "• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States",
    "• Built scalable web applications processing 50,000+ electoral boundaries with sub-200ms response times",
    "• Architected systems supporting 2,500+ concurrent users conducting redistricting analysis",
    "• Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"
)

# Create one new empty paragraph right after the target, then fill it (and the
# paragraphs implied by the embedded carriage returns) with the bullet text in
# a single InsertAfter call, so each bullet lands in its own <w:p>.
$target.Range.InsertParagraphAfter()
$insertPos = $target.Range.End
$d.Range($insertPos, $insertPos).InsertAfter([string]::Join("`r", $newBullets))
